$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells remain text (preserve formatting like "71.130.10")
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "71.130.10"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").Value = "3.858.07"
$ws.Range("E3").Value = "  -2.45%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "596.06"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").Value = "167.70"
$ws.Range("E6").Value = "  +6.96%  "
$ws.Range("D7").Value = "0.672"
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "0.751"
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").Value = "0.174"
$ws.Range("E10").Value = "  +4.38%  "
$ws.Range("D11").Value = "53.37"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "0.0000321"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "11.18"
$ws.Range("E13").Value = "  +3.55%  "
$ws.Range("D14").Value = "4.481.83"
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").Value = "3.874.34"
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("D16").Value = "20.71"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").Value = "13.79"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("E18").Value = "  -5.56%  "
$ws.Range("E19").Value = "  -2.04%  "
$ws.Range("D20").Value = "70.891.14"
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("D21").Value = "433.36"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").Value = "4.73"
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("D23").Value = "94.18"
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("E24").Value = "  -4.56%  "
$ws.Range("D25").Value = "13.77"
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("D26").Value = "4.11"
$ws.Range("E26").Value = "  -7.38%  "
$ws.Range("D27").Value = "10.91"
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("D28").Value = "5.93"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "10.18"
$ws.Range("E29").Value = "  -4.78%  "
$ws.Range("D30").Value = "34.96"
$ws.Range("E30").Value = "  -3.55%  "
$ws.Range("D31").Value = "7.94"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "13.49"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "49.60"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").Value = "0.125"
$ws.Range("E34").Value = "  -4.72%  "
$ws.Range("D35").Value = "69.10"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "0.0₃0979"
$ws.Range("E36").Value = "  +14.63%  "
$ws.Range("D37").Value = "622.86"
$ws.Range("E37").Value = "  -7.98%  "
$ws.Range("D38").Value = "0.418"
$ws.Range("E38").Value = "  -3.76%  "
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "3.31"
$ws.Range("E40").Value = "  +38.40%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "3.27"
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.142"
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("D44").Value = "0.0466"
$ws.Range("E44").Value = "  -3.74%  "
$ws.Range("D45").Value = "10.10"
$ws.Range("E45").Value = "  -6.90%  "
$ws.Range("D46").Value = "2.65"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D48").Value = "3.36"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").Value = "2.826.83"
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("D50").Value = "2.75"
$ws.Range("E50").Value = "  -18.00%  "
$ws.Range("E51").Value = "  +1.16%  "
